# Fix transmission and sets files
# Insert 20 new generator/node rows into the "GeneratorsOfNode" sheet,
# right before the existing "Moray Firth" / "Wind offshore grounded" row
# (previously row 718), shifting the subsequent rows down by 20
# (old A1:B753 -> new A1:B773).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneratorsOfNode")

# Insert 20 blank rows starting at row 718 (pushes old row 718.. down to 738..)
$ws.Range("A718:A737").Insert()

$newRows = @(
  @("Croatia", "Geo"),
  @("Bosnia H", "Hydro run-of-the-river"),
  @("Austria", "Nuclear"),
  @("Bosnia H", "Nuclear"),
  @("Denmark", "Nuclear"),
  @("Estonia", "Nuclear"),
  @("Greece", "Nuclear"),
  @("Croatia", "Nuclear"),
  @("Ireland", "Nuclear"),
  @("Italy", "Nuclear"),
  @("Luxemb.", "Nuclear"),
  @("Latvia", "Nuclear"),
  @("Macedonia", "Nuclear"),
  @("Portugal", "Nuclear"),
  @("Serbia", "Nuclear"),
  @("NO1", "Nuclear"),
  @("NO2", "Nuclear"),
  @("NO3", "Nuclear"),
  @("NO4", "Nuclear"),
  @("NO5", "Nuclear")
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 718 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
